$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Novas linhas de dados (fontes ajustadas para android/desktop)
$newRows = @(
    @(400, 400),
    @(500, 500),
    @(5000, 5000),
    @(8000, 8000),
    @(5000, 5000),
    @(5000, 8000)
)

$startRow = 33
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Ultima linha (39): valores numericos armazenados como texto ("5000.0")
$lastRow = 39
$ws.Range("A$lastRow").NumberFormat = "@"
$ws.Range("A$lastRow").Value = "5000.0"
$ws.Range("B$lastRow").NumberFormat = "@"
$ws.Range("B$lastRow").Value = "5000.0"
